$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly
# data (old D:K) right to F:M, to make room for two new quarters.
# Use a range bounded to the used rows (not EntireColumn) so the
# worksheet's used range / dimension doesn't balloon to 1048576 rows.
$ws.Range("D1:E102").Insert(-4161)

# The newly inserted D:E columns inherit formatting from column C by
# default; copy the number formats from the (now-shifted) F:G columns
# so D:E match the rest of the data columns (date format row 7/38/80,
# thousands-number format elsewhere).
$ws.Columns("F:G").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = most recent quarter,
# E = quarter before it) with the newly reported financial figures.
$ws.Cells.Item(7,4).Value2 = 43465
$ws.Cells.Item(7,5).Value2 = 43373
$ws.Cells.Item(8,4).Value2 = 29100
$ws.Cells.Item(8,5).Value2 = 30500
$ws.Cells.Item(9,4).Value2 = 26300
$ws.Cells.Item(9,5).Value2 = 28200
$ws.Cells.Item(10,4).Value2 = 2800
$ws.Cells.Item(10,5).Value2 = 2300
$ws.Cells.Item(12,4).Value2 = "NA"
$ws.Cells.Item(12,5).Value2 = "NA"
$ws.Cells.Item(13,4).Value2 = 0
$ws.Cells.Item(13,5).Value2 = 0
$ws.Cells.Item(14,4).Value2 = "NA"
$ws.Cells.Item(14,5).Value2 = "NA"
$ws.Cells.Item(15,4).Value2 = 400
$ws.Cells.Item(15,5).Value2 = 400
$ws.Cells.Item(17,4).Value2 = 30600
$ws.Cells.Item(17,5).Value2 = 32200
$ws.Cells.Item(18,4).Value2 = -1500
$ws.Cells.Item(18,5).Value2 = -1700
$ws.Cells.Item(20,4).Value2 = 0
$ws.Cells.Item(20,5).Value2 = 0
$ws.Cells.Item(21,4).Value2 = 400
$ws.Cells.Item(21,5).Value2 = 200
$ws.Cells.Item(22,4).Value2 = 300
$ws.Cells.Item(22,5).Value2 = 800
$ws.Cells.Item(23,4).Value2 = -1800
$ws.Cells.Item(23,5).Value2 = -2600
$ws.Cells.Item(24,4).Value2 = -500
$ws.Cells.Item(24,5).Value2 = 100
$ws.Cells.Item(25,4).Value2 = 0
$ws.Cells.Item(25,5).Value2 = 0
$ws.Cells.Item(26,4).Value2 = -1300
$ws.Cells.Item(26,5).Value2 = -2700
$ws.Cells.Item(27,4).Value2 = -1300
$ws.Cells.Item(27,5).Value2 = -2700
$ws.Cells.Item(28,4).Value2 = 0
$ws.Cells.Item(28,5).Value2 = 0
$ws.Cells.Item(29,4).Value2 = "NA"
$ws.Cells.Item(29,5).Value2 = 0
$ws.Cells.Item(30,4).Value2 = 0
$ws.Cells.Item(30,5).Value2 = 0
$ws.Cells.Item(31,4).Value2 = 0
$ws.Cells.Item(31,5).Value2 = 0
$ws.Cells.Item(32,4).Value2 = 0
$ws.Cells.Item(32,5).Value2 = 0
$ws.Cells.Item(33,4).Value2 = -1300
$ws.Cells.Item(33,5).Value2 = -2700
$ws.Cells.Item(34,4).Value2 = 0
$ws.Cells.Item(34,5).Value2 = 0
$ws.Cells.Item(35,4).Value2 = -1300
$ws.Cells.Item(35,5).Value2 = -2700
$ws.Cells.Item(38,4).Value2 = 43465
$ws.Cells.Item(38,5).Value2 = 43373
$ws.Cells.Item(41,4).Value2 = 900
$ws.Cells.Item(41,5).Value2 = 1300
$ws.Cells.Item(42,4).Value2 = 0
$ws.Cells.Item(42,5).Value2 = 0
$ws.Cells.Item(43,4).Value2 = 38000
$ws.Cells.Item(43,5).Value2 = 28100
$ws.Cells.Item(44,4).Value2 = 11800
$ws.Cells.Item(44,5).Value2 = 18300
$ws.Cells.Item(45,4).Value2 = 1500
$ws.Cells.Item(45,5).Value2 = 1900
$ws.Cells.Item(46,4).Value2 = 52200
$ws.Cells.Item(46,5).Value2 = 49600
$ws.Cells.Item(47,4).Value2 = 0
$ws.Cells.Item(47,5).Value2 = 0
$ws.Cells.Item(48,4).Value2 = 33800
$ws.Cells.Item(48,5).Value2 = 35400
$ws.Cells.Item(49,4).Value2 = 16400
$ws.Cells.Item(49,5).Value2 = 17100
$ws.Cells.Item(50,4).Value2 = 0
$ws.Cells.Item(50,5).Value2 = 0
$ws.Cells.Item(51,4).Value2 = 0
$ws.Cells.Item(51,5).Value2 = 0
$ws.Cells.Item(52,4).Value2 = 200
$ws.Cells.Item(52,5).Value2 = 200
$ws.Cells.Item(53,4).Value2 = 0
$ws.Cells.Item(53,5).Value2 = 0
$ws.Cells.Item(54,4).Value2 = 102600
$ws.Cells.Item(54,5).Value2 = 102200
$ws.Cells.Item(57,4).Value2 = 14200
$ws.Cells.Item(57,5).Value2 = 15500
$ws.Cells.Item(58,4).Value2 = 27200
$ws.Cells.Item(58,5).Value2 = 27200
$ws.Cells.Item(59,4).Value2 = 5200
$ws.Cells.Item(59,5).Value2 = 5100
$ws.Cells.Item(60,4).Value2 = 46600
$ws.Cells.Item(60,5).Value2 = 47800
$ws.Cells.Item(61,4).Value2 = 1900
$ws.Cells.Item(61,5).Value2 = 2300
$ws.Cells.Item(62,4).Value2 = 7600
$ws.Cells.Item(62,5).Value2 = 7900
$ws.Cells.Item(63,4).Value2 = 0
$ws.Cells.Item(63,5).Value2 = 0
$ws.Cells.Item(64,4).Value2 = 0
$ws.Cells.Item(64,5).Value2 = 0
$ws.Cells.Item(65,4).Value2 = 0
$ws.Cells.Item(65,5).Value2 = 0
$ws.Cells.Item(66,4).Value2 = 56200
$ws.Cells.Item(66,5).Value2 = 58000
$ws.Cells.Item(68,4).Value2 = 0
$ws.Cells.Item(68,5).Value2 = 0
$ws.Cells.Item(69,4).Value2 = 0
$ws.Cells.Item(69,5).Value2 = 0
$ws.Cells.Item(70,4).Value2 = 0
$ws.Cells.Item(70,5).Value2 = 0
$ws.Cells.Item(71,4).Value2 = 0
$ws.Cells.Item(71,5).Value2 = 0
$ws.Cells.Item(72,4).Value2 = 39400
$ws.Cells.Item(72,5).Value2 = 37100
$ws.Cells.Item(73,4).Value2 = 0
$ws.Cells.Item(73,5).Value2 = 0
$ws.Cells.Item(74,4).Value2 = 0
$ws.Cells.Item(74,5).Value2 = 0
$ws.Cells.Item(75,4).Value2 = 0
$ws.Cells.Item(75,5).Value2 = 0
$ws.Cells.Item(76,4).Value2 = 46400
$ws.Cells.Item(76,5).Value2 = 44200
$ws.Cells.Item(77,4).Value2 = 0
$ws.Cells.Item(77,5).Value2 = 0
$ws.Cells.Item(80,4).Value2 = 43465
$ws.Cells.Item(80,5).Value2 = 43373
$ws.Cells.Item(81,4).Value2 = -1300
$ws.Cells.Item(81,5).Value2 = -2700
$ws.Cells.Item(83,4).Value2 = 1900
$ws.Cells.Item(83,5).Value2 = 2000
$ws.Cells.Item(84,4).Value2 = 0
$ws.Cells.Item(84,5).Value2 = 0
$ws.Cells.Item(85,4).Value2 = 0
$ws.Cells.Item(85,5).Value2 = 0
$ws.Cells.Item(86,4).Value2 = 0
$ws.Cells.Item(86,5).Value2 = 0
$ws.Cells.Item(87,4).Value2 = 0
$ws.Cells.Item(87,5).Value2 = 0
$ws.Cells.Item(88,4).Value2 = 0
$ws.Cells.Item(88,5).Value2 = 0
$ws.Cells.Item(89,4).Value2 = 400
$ws.Cells.Item(89,5).Value2 = 600
$ws.Cells.Item(91,4).Value2 = -800
$ws.Cells.Item(91,5).Value2 = -1100
$ws.Cells.Item(92,4).Value2 = 0
$ws.Cells.Item(92,5).Value2 = 0
$ws.Cells.Item(93,4).Value2 = 0
$ws.Cells.Item(93,5).Value2 = 0
$ws.Cells.Item(94,4).Value2 = -500
$ws.Cells.Item(94,5).Value2 = -600
$ws.Cells.Item(96,4).Value2 = 0
$ws.Cells.Item(96,5).Value2 = 0
$ws.Cells.Item(97,4).Value2 = 0
$ws.Cells.Item(97,5).Value2 = 0
$ws.Cells.Item(98,4).Value2 = 0
$ws.Cells.Item(98,5).Value2 = 0
$ws.Cells.Item(99,4).Value2 = 0
$ws.Cells.Item(99,5).Value2 = 0
$ws.Cells.Item(100,4).Value2 = -300
$ws.Cells.Item(100,5).Value2 = -500
$ws.Cells.Item(101,4).Value2 = 0
$ws.Cells.Item(101,5).Value2 = 0
$ws.Cells.Item(102,4).Value2 = -300
$ws.Cells.Item(102,5).Value2 = -400

